$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values could otherwise be
# auto-converted to numbers by Excel (e.g. "1.00", "0.607", "6.61").
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "64.217.54"
$ws.Range("E2").Value = "  -2.73%  "
$ws.Range("D3").Value = "3.173.54"
$ws.Range("E3").Value = "  -7.92%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "562.97"
$ws.Range("E5").Value = "  -3.87%  "
$ws.Range("D6").Value = "170.43"
$ws.Range("E6").Value = "  -2.19%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "0.607"
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("D9").Value = "3.170.34"
$ws.Range("E9").Value = "  -7.99%  "
$ws.Range("E10").Value = "  -6.14%  "
$ws.Range("D11").Value = "6.61"
$ws.Range("E11").Value = "  -4.95%  "
$ws.Range("D12").Value = "0.396"
$ws.Range("E12").Value = "  -5.39%  "
$ws.Range("D13").Value = "3.724.21"
$ws.Range("E13").Value = "  -7.92%  "
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").Value = "27.36"
$ws.Range("E15").Value = "  -6.17%  "
$ws.Range("D16").Value = "64.245.80"
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("E17").Value = "  -5.57%  "
$ws.Range("D18").Value = "3.177.22"
$ws.Range("E18").Value = "  -7.79%  "
$ws.Range("D19").Value = "5.69"
$ws.Range("E19").Value = "  -4.59%  "
$ws.Range("D20").Value = "13.05"
$ws.Range("E20").Value = "  -5.77%  "
$ws.Range("D21").Value = "352.65"
$ws.Range("E21").Value = "  -4.95%  "
$ws.Range("D22").Value = "7.19"
$ws.Range("E22").Value = "  -5.38%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "69.07"
$ws.Range("E24").Value = "  -4.65%  "
$ws.Range("D25").Value = "0.502"
$ws.Range("E25").Value = "  -5.75%  "
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("D27").Value = "9.58"
$ws.Range("E27").Value = "  -1.69%  "
$ws.Range("E28").Value = "  -1.95%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").Value = "5.64"
$ws.Range("E30").Value = "  -2.79%  "
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").Value = "1.90"
$ws.Range("E32").Value = "  -4.46%  "
$ws.Range("D33").Value = "22.11"
$ws.Range("E33").Value = "  -6.55%  "
$ws.Range("D34").Value = "6.65"
$ws.Range("E34").Value = "  -5.43%  "
$ws.Range("D35").Value = "1.20"
$ws.Range("E35").Value = "  -5.54%  "
$ws.Range("E36").Value = "  -6.14%  "
$ws.Range("D37").Value = "155.48"
$ws.Range("E37").Value = "  -3.81%  "
$ws.Range("D38").Value = "0.810"
$ws.Range("E38").Value = "  -8.05%  "
$ws.Range("D39").Value = "26.02"
$ws.Range("E39").Value = "  -8.07%  "
$ws.Range("E40").Value = "  -2.56%  "
$ws.Range("E41").Value = "  -4.27%  "
$ws.Range("D42").Value = "2.608.08"
$ws.Range("E42").Value = "  -6.58%  "
$ws.Range("D43").Value = "4.18"
$ws.Range("E43").Value = "  -6.49%  "
$ws.Range("D44").Value = "6.04"
$ws.Range("E44").Value = "  -6.74%  "
$ws.Range("D45").Value = "331.55"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("D46").Value = "0.0655"
$ws.Range("E46").Value = "  -4.75%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "23.99"
$ws.Range("E47").Value = "  -4.71%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "38.84"
$ws.Range("E48").Value = "  -2.39%  "
$ws.Range("D49").Value = "0.0271"
$ws.Range("E49").Value = "  -7.23%  "
$ws.Range("D50").Value = "0.101"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("E51").Value = "  -0.03%  "
